# LITE-23297: sync translation attributes on product sync
#
# The workbook currently has a single per-language translation sheet
# ("TRN-1079-0833-9890 (FA)"). We duplicate it to add a second
# translation language sheet ("ES (TRN-1079-0833-9891)") with the same
# attribute rows/layout, then rename both sheets so the TRN id moves to
# the end of the tab name, and make the newly added sheet the active one.

$wb = $excel.ActiveWorkbook

$faSheet = $wb.Worksheets.Item("TRN-1079-0833-9890 (FA)")

# Duplicate the FA translation sheet right after itself - this becomes
# the new ES translation sheet.
$faSheet.Copy($null, $faSheet)
$esSheet = $wb.Worksheets.Item($faSheet.Index + 1)

# Rename both sheets: TRN id moves from prefix to parenthesized suffix.
$faSheet.Name = "FA (TRN-1079-0833-9890)"
$esSheet.Name = "ES (TRN-1079-0833-9891)"

# Re-create the hidden _FilterDatabase defined name scoped to the new
# sheet (mirrors the one already present on the FA sheet for its
# autofilter range).
$esSheet.Names.Add("_xlnm._FilterDatabase", "='ES (TRN-1079-0833-9891)'!`$A:`$E", $false)

# Make the new ES sheet the active / selected tab.
$esSheet.Activate()
